$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 386 ("象の群れ。ドローンの音から逃げている") entirely.
# This shifts all subsequent rows up by one (old row 387 becomes new row 386, etc.)
$ws.Rows.Item(386).Delete()
